$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "{'C': 100, 'degree': 2, 'gamma': 'scale', 'kernel': 'rbf'}"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "65.75%"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "66.56%"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "65.75%"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "65.19%"
